# Antes de uma arrumacao geral
#
# - Remove the (empty/unused) title slide that used to be the first slide
#   in the deck, leaving the diagram slide as the only slide.
# - Nudge two rectangles inside the "biblatex vs biber/natbib" diagram
#   group on the remaining slide so their boxes better fit the text they
#   contain (the "biblatex" box gets a touch narrower, the "natbib" box
#   slides left and gets wider).

$p = $ppt.ActivePresentation

# 1) Delete the original first slide (title/subtitle placeholder slide).
#    The deck's only remaining slide (formerly slide 2) becomes slide 1.
$p.Slides.Item(1).Delete()

# 2) Resize/reposition two shapes that live inside the "Group 26" group on
#    the (now only) slide.
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item("Group 26")

# "Rectangle 1" ("biblatex" box): keep its position, shrink its width
# from 4146900 EMU to 4097832 EMU.
$rect1 = $grp.GroupItems.Item("Rectangle 1")
$rect1.Width = 322.663970957874

# "Rectangle 5" ("natbib" box): shift left (off x 7239700 -> 7110484 EMU)
# and widen it (cx 1693951 -> 1823168 EMU).
$rect5 = $grp.GroupItems.Item("Rectangle 5")
$rect5.Left = 559.8806762712599
$rect5.Width = 143.55657963307087
